$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Case_4_20 with 380 kV slack bus voltage set-point lowered from 1.05 pu to 1.02 pu.
# Column B holds the slack-bus set-point per row; columns C-F and I-N hold the
# recalculated per-unit voltage magnitudes at the other buses for that scenario.
$data = @{
    2 = @{ "B" = 1.02; "C" = 1.031483913438708; "D" = 1.037752324847416; "E" = 1.035122920090494; "F" = 1.045208720303641; "I" = 1.032424074483189; "J" = 1.036619387477876; "K" = 1.040542312718022; "L" = 1.037920428664011; "M" = 1.047977604476914; "N" = 1.038091505910029 }
    3 = @{ "B" = 1.02; "C" = 1.032556738882705; "D" = 1.038558565408314; "E" = 1.036143898266447; "F" = 1.046218500185343; "I" = 1.032590958305202; "J" = 1.037333479179003; "K" = 1.04115843521033; "L" = 1.038750173503523; "M" = 1.048798261628608; "N" = 1.03880661170325 }
    4 = @{ "B" = 1.02; "C" = 1.033251204644971; "D" = 1.03908028757045; "E" = 1.03680517364483; "F" = 1.046872404035093; "I" = 1.032697566401658; "J" = 1.037795269105075; "K" = 1.041556487413829; "L" = 1.039287103845; "M" = 1.049329180945313; "N" = 1.039269057423986 }
    5 = @{ "B" = 1.02; "C" = 1.033543224240139; "D" = 1.039299626172374; "E" = 1.037083325269095; "F" = 1.047147426449089; "I" = 1.032742054505131; "J" = 1.037989339408449; "K" = 1.041723679594103; "L" = 1.039512836193471; "M" = 1.049552354811513; "N" = 1.03946340332945 }
    6 = @{ "B" = 1.02; "C" = 1.03359225950014; "D" = 1.039336454455085; "E" = 1.037130037021148; "F" = 1.047193611015668; "I" = 1.032749504896648; "J" = 1.038021920776474; "K" = 1.041751743147545; "L" = 1.039550738039003; "M" = 1.049589825227807; "N" = 1.039496030966755 }
    7 = @{ "B" = 1.02; "C" = 1.033255106364853; "D" = 1.039083218358543; "E" = 1.036808889725806; "F" = 1.046876078422377; "I" = 1.032698162150699; "J" = 1.037797862542404; "K" = 1.041558722029863; "L" = 1.039290120065261; "M" = 1.049332163102781; "N" = 1.039271654544292 }
    8 = @{ "B" = 1.02; "C" = 1.031846422753253; "D" = 1.038024790847954; "E" = 1.035467832895003; "F" = 1.045549874466909; "I" = 1.03248075856169; "J" = 1.036860774897966; "K" = 1.040750662338113; "L" = 1.038200838413872; "M" = 1.048254969961645; "N" = 1.038333236127914 }
    9 = @{ "B" = 1.02; "C" = 1.029366249937053; "D" = 1.036159972341695; "E" = 1.033109593279611; "F" = 1.04321685435393; "I" = 1.032087133808021; "J" = 1.035207412719884; "K" = 1.039322028357345; "L" = 1.036281628854974; "M" = 1.046356068165256; "N" = 1.036677525985896 }
    10 = @{ "B" = 1.02; "C" = 1.027714207454214; "D" = 1.034916979659439; "E" = 1.031540730197007; "F" = 1.04166417420911; "I" = 1.031817653815993; "J" = 1.034103779398431; "K" = 1.038366452205818; "L" = 1.035002337242652; "M" = 1.045089656679003; "N" = 1.035572325378538 }
    11 = @{ "B" = 1.02; "C" = 1.026999184959289; "D" = 1.03437881023551; "E" = 1.030862180466008; "F" = 1.040992484748903; "I" = 1.031699294091008; "J" = 1.033625566591926; "K" = 1.037951934234344; "L" = 1.034448435863006; "M" = 1.044541178472536; "N" = 1.035093433454996 }
    12 = @{ "B" = 1.02; "C" = 1.026733641868928; "D" = 1.034178918931572; "E" = 1.030610254001536; "F" = 1.040743084594117; "I" = 1.031655078886125; "J" = 1.033447887278191; "K" = 1.037797851877177; "L" = 1.034242698441227; "M" = 1.044337432493629; "N" = 1.034915501816255 }
    13 = @{ "B" = 1.02; "C" = 1.026790599592891; "D" = 1.034221795888143; "E" = 1.030664287780579; "F" = 1.040796577459117; "I" = 1.031664574558219; "J" = 1.033486002344771; "K" = 1.037830908125412; "L" = 1.034286829532952; "M" = 1.044381137450552; "N" = 1.034953671010602 }
    14 = @{ "B" = 1.02; "C" = 1.026977234105642; "D" = 1.034362286975481; "E" = 1.030841353755635; "F" = 1.040971867309729; "I" = 1.031695644370067; "J" = 1.033610880579217; "K" = 1.037939200019978; "L" = 1.034431429415902; "M" = 1.044524337106205; "N" = 1.035078726586464 }
    15 = @{ "B" = 1.02; "C" = 1.027092232150676; "D" = 1.034448849352154; "E" = 1.030950465471527; "F" = 1.041079881800455; "I" = 1.031714754235992; "J" = 1.033687815588138; "K" = 1.038005907408909; "L" = 1.034520523004107; "M" = 1.044612564930212; "N" = 1.035155770851923 }
    16 = @{ "B" = 1.02; "C" = 1.027761667909442; "D" = 1.034952697391583; "E" = 1.031585779722045; "F" = 1.041708765358837; "I" = 1.031825473730481; "J" = 1.034135509818172; "K" = 1.038393946684178; "L" = 1.035039098715813; "M" = 1.045126055004031; "N" = 1.035604100859115 }
    17 = @{ "B" = 1.02; "C" = 1.028181673416825; "D" = 1.035268762984103; "E" = 1.031984504077494; "F" = 1.042103416677783; "I" = 1.031894477369303; "J" = 1.034416247713545; "K" = 1.038637153658946; "L" = 1.035364398465218; "M" = 1.045448123557296; "N" = 1.0358852374345 }
    18 = @{ "B" = 1.02; "C" = 1.028426686695461; "D" = 1.035453123927221; "E" = 1.032217148334096; "F" = 1.042333671000257; "I" = 1.031934564577858; "J" = 1.034579965307164; "K" = 1.038778940057192; "L" = 1.035554144186342; "M" = 1.045635969705449; "N" = 1.036049187525881 }
    19 = @{ "B" = 1.02; "C" = 1.028510235191219; "D" = 1.035515987114089; "E" = 1.032296486743898; "F" = 1.042412192089724; "I" = 1.031948205888825; "J" = 1.034635783331737; "K" = 1.038827273312814; "L" = 1.035618843227385; "M" = 1.045700018513571; "N" = 1.036105084818451 }
    20 = @{ "B" = 1.02; "C" = 1.02813660757923; "D" = 1.035234851591987; "E" = 1.031941716959738; "F" = 1.042061068025776; "I" = 1.031887090629362; "J" = 1.034386130510421; "K" = 1.038611067309012; "L" = 1.035329496478087; "M" = 1.04541356978386; "N" = 1.035855077461494 }
    21 = @{ "B" = 1.02; "C" = 1.026922273568732; "D" = 1.034320915604846; "E" = 1.030789209004735; "F" = 1.040920246186228; "I" = 1.0316865020155; "J" = 1.03357410842487; "K" = 1.037907313843183; "L" = 1.034388848181676; "M" = 1.044482168836585; "N" = 1.035041902211441 }
    22 = @{ "B" = 1.02; "C" = 1.026159051681693; "D" = 1.033746338871406; "E" = 1.030065259210555; "F" = 1.040203517302574; "I" = 1.031558931130023; "J" = 1.033063269939591; "K" = 1.037464188701148; "L" = 1.033797461287626; "M" = 1.043896463333898; "N" = 1.034530338276919 }
    23 = @{ "B" = 1.02; "C" = 1.026563623815226; "D" = 1.034050927833344; "E" = 1.030448974394865; "F" = 1.040583416483215; "I" = 1.031626696506802; "J" = 1.033334102335811; "K" = 1.037699158979246; "L" = 1.034110963203233; "M" = 1.044206965965366; "N" = 1.034801555286206 }
    24 = @{ "B" = 1.02; "C" = 1.028156970794874; "D" = 1.035250174674204; "E" = 1.031961050384399; "F" = 1.042080203370547; "I" = 1.031890428878014; "J" = 1.034399739277139; "K" = 1.038622854830447; "L" = 1.035345267171647; "M" = 1.045429183180615; "N" = 1.035868705554221 }
    25 = @{ "B" = 1.02; "C" = 1.03000718543574; "D" = 1.036642036891644; "E" = 1.033718675137204; "F" = 1.043819528206198; "I" = 1.032190141515703; "J" = 1.035635092975686; "K" = 1.039691922241914; "L" = 1.036777759831591; "M" = 1.046847065854627; "N" = 1.037105813596717 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
